$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(3409, 45565, 28192.01, 45.66, 5.89, 0.68),
  @(3410, 45566, 28277.68, 45.8, 5.91, 0.68),
  @(3411, 45568, 27677.44, 44.83, 5.78, 0.7),
  @(3412, 45569, 27432.2, 44.43, 5.73, 0.7),
  @(3413, 45572, 26868.49, 43.52, 5.61, 0.72),
  @(3414, 45573, 27400.64, 44.38, 5.72, 0.7),
  @(3415, 45574, 27694.62, 44.85, 5.78, 0.7),
  @(3416, 45575, 27653.34, 44.49, 5.71, 0.7),
  @(3417, 45576, 27772.52, 44.69, 5.73, 0.7),
  @(3418, 45579, 27873.16, 44.85, 5.75, 0.6899999999999999),
  @(3419, 45580, 27961.43, 44.96, 5.77, 0.6899999999999999),
  @(3420, 45581, 27906.39, 44.88, 5.76, 0.6899999999999999),
  @(3421, 45582, 27443.15, 44.06, 5.67, 0.7),
  @(3422, 45583, 27490.06, 44.13, 5.68, 0.7),
  @(3423, 45586, 27078.49, 43.56, 5.59, 0.71),
  @(3424, 45587, 26382.91, 42.5, 5.45, 0.73),
  @(3425, 45588, 26516.6, 42.72, 5.47, 0.73),
  @(3426, 45589, 26450.31, 42.03, 5.46, 0.73),
  @(3427, 45590, 25984.56, 41.11, 5.36, 0.75),
  @(3428, 45593, 26151.21, 41.23, 5.4, 0.74),
  @(3429, 45594, 26368.32, 43.39, 5.44, 0.73),
  @(3430, 45595, 26432.02, 43.5, 5.46, 0.73),
  @(3431, 45596, 26381.69, 42.98, 5.44, 0.74),
  @(3432, 45597, 26564.21, 43.28, 5.48, 0.73),
  @(3433, 45600, 26246.88, 42.88, 5.41, 0.75),
  @(3434, 45601, 26366.59, 43.06, 5.44, 0.75),
  @(3435, 45602, 26920.14, 44, 5.55, 0.73),
  @(3436, 45603, 26789.57, 43.45, 5.52, 0.73),
  @(3437, 45604, 26442.2, 42.87, 5.45, 0.75),
  @(3438, 45607, 26227.35, 42.31, 5.41, 0.74),
  @(3439, 45608, 25989.54, 41.9, 5.36, 0.75),
  @(3440, 45609, 25294.53, 40.73, 5.21, 0.78),
  @(3441, 45610, 25417.32, 40.9, 5.24, 0.78),
  @(3442, 45614, 25371.84, 40.64, 5.23, 0.77),
  @(3443, 45615, 25614.67, 40.69, 5.28, 0.78),
  @(3444, 45617, 25517.07, 40.54, 5.26, 0.77),
  @(3445, 45618, 25810.08, 41, 5.32, 0.76),
  @(3446, 45621, 26205.9, 41.63, 5.4, 0.75),
  @(3447, 45622, 26218.11, 41.64, 5.4, 0.78),
  @(3448, 45623, 26348.67, 41.84, 5.37, 0.78),
  @(3449, 45624, 26351.73, 41.85, 5.37, 0.78),
  @(3450, 45625, 26436.15, 41.98, 5.38, 0.78),
  @(3451, 45628, 26688.15, 42.38, 5.44, 0.77),
  @(3452, 45629, 26922.62, 42.75, 5.48, 0.77),
  @(3453, 45630, 27167.68, 43.14, 5.53, 0.76),
  @(3454, 45631, 27284.64, 43.33, 5.56, 0.76),
  @(3455, 45632, 27395.59, 43.51, 5.58, 0.75),
  @(3456, 45635, 27500.64, 43.67, 5.6, 0.75),
  @(3457, 45636, 27577.27, 43.79, 5.62, 0.75),
  @(3458, 45637, 27666.35, 43.94, 5.64, 0.75),
  @(3459, 45638, 27527.93, 43.72, 5.61, 0.75),
  @(3460, 45639, 27500.16, 43.67, 5.6, 0.75),
  @(3461, 45642, 27690.73, 43.97, 5.64, 0.74),
  @(3462, 45643, 27507.72, 43.68, 5.6, 0.75),
  @(3463, 45644, 27342.13, 43.42, 5.57, 0.75),
  @(3464, 45645, 27276.63, 43.32, 5.56, 0.75),
  @(3465, 45646, 26618.21, 42.27, 5.42, 0.77),
  @(3466, 45649, 26650.5, 42.32, 5.43, 0.77),
  @(3467, 45650, 26656.27, 42.33, 5.43, 0.77),
  @(3468, 45652, 26676.67, 42.36, 5.43, 0.77),
  @(3469, 45653, 26645.08, 42.31, 5.43, 0.77),
  @(3470, 45656, 26692.33, 42.39, 5.44, 0.77),
  @(3471, 45657, 26732.77, 43.3, 5.53, 0.76),
  @(3472, 45658, 26865.92, 43.52, 5.56, 0.76),
  @(3473, 45659, 27152.47, 43.98, 5.62, 0.75),
  @(3474, 45660, 27084.39, 43.87, 5.6, 0.75),
  @(3475, 45663, 26402.02, 42.77, 5.46, 0.77),
  @(3476, 45664, 26631.09, 43.14, 5.51, 0.76),
  @(3477, 45665, 26352.27, 42.69, 5.45, 0.77),
  @(3478, 45666, 26142.28, 42.35, 5.41, 0.78),
  @(3479, 45667, 25609.55, 41.48, 5.3, 0.79),
  @(3480, 45670, 24604.76, 39.86, 5.09, 0.83),
  @(3481, 45671, 25142.23, 40.72, 5.2, 0.8100000000000001)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
    $ws.Cells.Item($r, 5).Value2 = $row[5]
}

Write-Host "Added" $data.Count "rows"
